$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Add new "time_taken" header in F1 (matching style of other headers)
# ---------------------------------------------------------------
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------
# 2. Full gene-row data (columns B,C,D,E,F) for rows 2..226 (225 rows)
#    Row order accounts for two newly-inserted genes:
#      - RFXANK  (inserted so it ends up at row 187)
#      - DAB1    (inserted so it ends up at row 209)
#    A leading "'" on the geneConfidence (D) field forces it to be
#    stored as text (e.g. "1", "2", "3"), matching the source data.
# ---------------------------------------------------------------
$text = @'
AAAS|aladin WD repeat nucleoporin|'3|Ataxia - paediatric|2021-10-05 10:50:11.459087
ABCB7|ATP binding cassette subfamily B member 7|'3|Ataxia - paediatric|2021-10-05 10:50:11.459098
ABHD12|abhydrolase domain containing 12|'3|Ataxia - paediatric|2021-10-05 10:50:11.459102
ACO2|aconitase 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459104
ADGRG1|adhesion G protein-coupled receptor G1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459107
ADPRHL2|ADP-ribosylhydrolase like 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459110
AFG3L2|AFG3 like matrix AAA peptidase subunit 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459112
AHI1|Abelson helper integration site 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459115
ALDH5A1|aldehyde dehydrogenase 5 family member A1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459118
ANO10|anoctamin 10|'3|Ataxia - paediatric|2021-10-05 10:50:11.459120
AP1S2|adaptor related protein complex 1 sigma 2 subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459122
APTX|aprataxin|'3|Ataxia - paediatric|2021-10-05 10:50:11.459125
ARL13B|ADP ribosylation factor like GTPase 13B|'3|Ataxia - paediatric|2021-10-05 10:50:11.459127
ARSA|arylsulfatase A|'3|Ataxia - paediatric|2021-10-05 10:50:11.459130
ATCAY|ATCAY, caytaxin|'3|Ataxia - paediatric|2021-10-05 10:50:11.459132
ATG7|autophagy related 7|'3|Ataxia - paediatric|2021-10-05 10:50:11.459135
ATM|ATM serine/threonine kinase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459138
ATP1A3|ATPase Na+/K+ transporting subunit alpha 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459140
ATP8A2|ATPase phospholipid transporting 8A2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459143
BBS1|Bardet-Biedl syndrome 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459145
C5orf42|chromosome 5 open reading frame 42|'3|Ataxia - paediatric|2021-10-05 10:50:11.459147
CA8|carbonic anhydrase 8|'3|Ataxia - paediatric|2021-10-05 10:50:11.459150
CACNA1A|calcium voltage-gated channel subunit alpha1 A|'3|Ataxia - paediatric|2021-10-05 10:50:11.459152
CACNA1G|calcium voltage-gated channel subunit alpha1 G|'3|Ataxia - paediatric|2021-10-05 10:50:11.459155
CAD|carbamoyl-phosphate synthetase 2, aspartate transcarbamylase, and dihydroorotase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459158
CAMTA1|calmodulin binding transcription activator 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459160
CASK|calcium/calmodulin dependent serine protein kinase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459163
CBY1|chibby family member 1, beta catenin antagonist|'3|Ataxia - paediatric|2021-10-05 10:50:11.459165
CC2D2A|coiled-coil and C2 domain containing 2A|'3|Ataxia - paediatric|2021-10-05 10:50:11.459168
CEP290|centrosomal protein 290|'3|Ataxia - paediatric|2021-10-05 10:50:11.459170
CEP41|centrosomal protein 41|'3|Ataxia - paediatric|2021-10-05 10:50:11.459172
CLN5|CLN5, intracellular trafficking protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459175
CLN6|CLN6, transmembrane ER protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459178
CLPP|caseinolytic mitochondrial matrix peptidase proteolytic subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459180
COA7|cytochrome c oxidase assembly factor 7 (putative)|'3|Ataxia - paediatric|2021-10-05 10:50:11.459182
COQ8A|coenzyme Q8A|'3|Ataxia - paediatric|2021-10-05 10:50:11.459185
COX20|COX20, cytochrome c oxidase assembly factor|'3|Ataxia - paediatric|2021-10-05 10:50:11.459187
CSTB|cystatin B|'3|Ataxia - paediatric|2021-10-05 10:50:11.459190
CTBP1|C-terminal binding protein 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459192
CWF19L1|CWF19 like 1, cell cycle control (S. pombe)|'3|Ataxia - paediatric|2021-10-05 10:50:11.459195
CYP27A1|cytochrome P450 family 27 subfamily A member 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459197
DARS2|aspartyl-tRNA synthetase 2, mitochondrial|'3|Ataxia - paediatric|2021-10-05 10:50:11.459200
DDHD2|DDHD domain containing 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459202
DNAJC19|DnaJ heat shock protein family (Hsp40) member C19|'3|Ataxia - paediatric|2021-10-05 10:50:11.459205
DOCK3|dedicator of cytokinesis 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459207
EBF3|early B-cell factor 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459210
EIF2B1|eukaryotic translation initiation factor 2B subunit alpha|'3|Ataxia - paediatric|2021-10-05 10:50:11.459212
EIF2B2|eukaryotic translation initiation factor 2B subunit beta|'3|Ataxia - paediatric|2021-10-05 10:50:11.459215
EIF2B3|eukaryotic translation initiation factor 2B subunit gamma|'3|Ataxia - paediatric|2021-10-05 10:50:11.459217
EIF2B4|eukaryotic translation initiation factor 2B subunit delta|'3|Ataxia - paediatric|2021-10-05 10:50:11.459220
EIF2B5|eukaryotic translation initiation factor 2B subunit epsilon|'3|Ataxia - paediatric|2021-10-05 10:50:11.459222
EPM2A|EPM2A, laforin glucan phosphatase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459225
EXOSC5|exosome component 5|'3|Ataxia - paediatric|2021-10-05 10:50:11.459228
FA2H|fatty acid 2-hydroxylase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459230
FBXL4|F-box and leucine rich repeat protein 4|'3|Ataxia - paediatric|2021-10-05 10:50:11.459233
FGF14|fibroblast growth factor 14|'3|Ataxia - paediatric|2021-10-05 10:50:11.459235
FLVCR1|feline leukemia virus subgroup C cellular receptor 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459238
FOLR1|folate receptor 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459240
FXN|frataxin|'3|Ataxia - paediatric|2021-10-05 10:50:11.459242
GBA2|glucosylceramidase beta 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459245
GJC2|gap junction protein gamma 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459247
GOSR2|golgi SNAP receptor complex member 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459250
GPAA1|glycosylphosphatidylinositol anchor attachment 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459252
GRID2|glutamate ionotropic receptor delta type subunit 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459255
GRM1|glutamate metabotropic receptor 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459258
GSS|glutathione synthetase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459261
HEXA|hexosaminidase subunit alpha|'3|Ataxia - paediatric|2021-10-05 10:50:11.459263
HEXB|hexosaminidase subunit beta|'3|Ataxia - paediatric|2021-10-05 10:50:11.459266
INPP5E|inositol polyphosphate-5-phosphatase E|'3|Ataxia - paediatric|2021-10-05 10:50:11.459268
IRF2BPL|interferon regulatory factor 2 binding protein like|'3|Ataxia - paediatric|2021-10-05 10:50:11.459270
ITPR1|inositol 1,4,5-trisphosphate receptor type 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459273
KCNA1|potassium voltage-gated channel subfamily A member 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459275
KCNA2|potassium voltage-gated channel subfamily A member 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459278
KCNC3|potassium voltage-gated channel subfamily C member 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459280
KCNJ10|potassium voltage-gated channel subfamily J member 10|'3|Ataxia - paediatric|2021-10-05 10:50:11.459283
KCNN2|potassium calcium-activated channel subfamily N member 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459285
KIF1C|kinesin family member 1C|'3|Ataxia - paediatric|2021-10-05 10:50:11.459289
KIF7|kinesin family member 7|'3|Ataxia - paediatric|2021-10-05 10:50:11.459292
LAMA1|laminin subunit alpha 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459295
LARS2|leucyl-tRNA synthetase 2, mitochondrial|'3|Ataxia - paediatric|2021-10-05 10:50:11.459297
MAG|myelin associated glycoprotein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459300
MARS2|methionyl-tRNA synthetase 2, mitochondrial|'3|Ataxia - paediatric|2021-10-05 10:50:11.459302
MKS1|Meckel syndrome, type 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459305
MMACHC|methylmalonic aciduria (cobalamin deficiency) cblC type, with homocystinuria|'3|Ataxia - paediatric|2021-10-05 10:50:11.459307
MORC2|MORC family CW-type zinc finger 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459309
MRE11|MRE11 homolog, double strand break repair nuclease|'3|Ataxia - paediatric|2021-10-05 10:50:11.459312
MSTO1|misato 1, mitochondrial distribution and morphology regulator|'3|Ataxia - paediatric|2021-10-05 10:50:11.459314
MTCL1|microtubule crosslinking factor 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459317
MTFMT|mitochondrial methionyl-tRNA formyltransferase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459319
MTTP|microsomal triglyceride transfer protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459322
MVK|mevalonate kinase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459324
NHLRC1|NHL repeat containing E3 ubiquitin protein ligase 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459326
NKX2-1|NK2 homeobox 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459330
NKX6-2|NK6 homeobox 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459333
NPC1|NPC intracellular cholesterol transporter 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459335
NPC2|NPC intracellular cholesterol transporter 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459338
NPHP1|nephrocystin 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459340
NUS1|NUS1 dehydrodolichyl diphosphate synthase subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459343
OFD1|OFD1, centriole and centriolar satellite protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459345
OPA1|OPA1, mitochondrial dynamin like GTPase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459347
OPA3|OPA3, outer mitochondrial membrane lipid metabolism regulator|'3|Ataxia - paediatric|2021-10-05 10:50:11.459350
OPHN1|oligophrenin 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459352
PAX6|paired box 6|'3|Ataxia - paediatric|2021-10-05 10:50:11.459355
PEX16|peroxisomal biogenesis factor 16|'3|Ataxia - paediatric|2021-10-05 10:50:11.459357
PEX7|peroxisomal biogenesis factor 7|'3|Ataxia - paediatric|2021-10-05 10:50:11.459360
PHYH|phytanoyl-CoA 2-hydroxylase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459362
PITRM1|pitrilysin metallopeptidase 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459365
PLA2G6|phospholipase A2 group VI|'3|Ataxia - paediatric|2021-10-05 10:50:11.459367
PMPCA|peptidase, mitochondrial processing alpha subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459371
PMPCB|peptidase, mitochondrial processing beta subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459374
PNKD|paroxysmal nonkinesigenic dyskinesia|'3|Ataxia - paediatric|2021-10-05 10:50:11.459377
PNKP|polynucleotide kinase 3'-phosphatase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459379
PNPLA6|patatin like phospholipase domain containing 6|'3|Ataxia - paediatric|2021-10-05 10:50:11.459382
POLG|DNA polymerase gamma, catalytic subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459384
POLR3A|RNA polymerase III subunit A|'3|Ataxia - paediatric|2021-10-05 10:50:11.459387
POLR3B|RNA polymerase III subunit B|'3|Ataxia - paediatric|2021-10-05 10:50:11.459389
POU4F1|POU class 4 homeobox 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459392
PRDX3|peroxiredoxin 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459394
PRICKLE1|prickle planar cell polarity protein 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459397
PRRT2|proline rich transmembrane protein 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459399
PTRH2|peptidyl-tRNA hydrolase 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459402
RNF216|ring finger protein 216|'3|Ataxia - paediatric|2021-10-05 10:50:11.459404
RNF220|ring finger protein 220|'3|Ataxia - paediatric|2021-10-05 10:50:11.459407
RORA|RAR related orphan receptor A|'3|Ataxia - paediatric|2021-10-05 10:50:11.459409
RPGRIP1L|RPGRIP1 like|'3|Ataxia - paediatric|2021-10-05 10:50:11.459411
RUBCN|RUN and cysteine rich domain containing beclin 1 interacting protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459414
SACS|sacsin molecular chaperone|'3|Ataxia - paediatric|2021-10-05 10:50:11.459416
SCN1A|sodium voltage-gated channel alpha subunit 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459419
SCN2A|sodium voltage-gated channel alpha subunit 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459423
SCN8A|sodium voltage-gated channel alpha subunit 8|'3|Ataxia - paediatric|2021-10-05 10:50:11.459426
SCYL1|SCY1 like pseudokinase 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459429
SETX|senataxin|'3|Ataxia - paediatric|2021-10-05 10:50:11.459431
SIL1|SIL1 nucleotide exchange factor|'3|Ataxia - paediatric|2021-10-05 10:50:11.459433
SLC17A5|solute carrier family 17 member 5|'3|Ataxia - paediatric|2021-10-05 10:50:11.459436
SLC1A3|solute carrier family 1 member 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459438
SLC25A46|solute carrier family 25 member 46|'3|Ataxia - paediatric|2021-10-05 10:50:11.459441
SLC2A1|solute carrier family 2 member 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459443
SLC44A1|solute carrier family 44 member 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459446
SLC52A2|solute carrier family 52 member 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459448
SLC9A6|solute carrier family 9 member A6|'3|Ataxia - paediatric|2021-10-05 10:50:11.459450
SNAP25|synaptosome associated protein 25|'3|Ataxia - paediatric|2021-10-05 10:50:11.459453
SNX14|sorting nexin 14|'3|Ataxia - paediatric|2021-10-05 10:50:11.459455
SPR|sepiapterin reductase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459458
SPTBN2|spectrin beta, non-erythrocytic 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459460
SQSTM1|sequestosome 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459463
SRD5A3|steroid 5 alpha-reductase 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459465
STUB1|STIP1 homology and U-box containing protein 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459468
SVBP|small vasohibin binding protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459470
SYNE1|spectrin repeat containing nuclear envelope protein 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459473
TBC1D23|TBC1 domain family member 23|'3|Ataxia - paediatric|2021-10-05 10:50:11.459476
TCTN1|tectonic family member 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459478
TCTN2|tectonic family member 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459481
TCTN3|tectonic family member 3|'3|Ataxia - paediatric|2021-10-05 10:50:11.459483
TDP2|tyrosyl-DNA phosphodiesterase 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459486
THG1L|tRNA-histidine guanylyltransferase 1 like|'3|Ataxia - paediatric|2021-10-05 10:50:11.459488
TINF2|TERF1 interacting nuclear factor 2|'3|Ataxia - paediatric|2021-10-05 10:50:11.459491
TMEM106B|transmembrane protein 106B|'3|Ataxia - paediatric|2021-10-05 10:50:11.459493
TMEM216|transmembrane protein 216|'3|Ataxia - paediatric|2021-10-05 10:50:11.459496
TMEM237|transmembrane protein 237|'3|Ataxia - paediatric|2021-10-05 10:50:11.459498
TMEM240|transmembrane protein 240|'3|Ataxia - paediatric|2021-10-05 10:50:11.459501
TMEM67|transmembrane protein 67|'3|Ataxia - paediatric|2021-10-05 10:50:11.459503
TPP1|tripeptidyl peptidase 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459506
TSFM|Ts translation elongation factor, mitochondrial|'3|Ataxia - paediatric|2021-10-05 10:50:11.459508
TTC19|tetratricopeptide repeat domain 19|'3|Ataxia - paediatric|2021-10-05 10:50:11.459511
TTPA|alpha tocopherol transfer protein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459513
TUBB4A|tubulin beta 4A class IVa|'3|Ataxia - paediatric|2021-10-05 10:50:11.459516
TWNK|twinkle mtDNA helicase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459518
UBTF|upstream binding transcription factor, RNA polymerase I|'3|Ataxia - paediatric|2021-10-05 10:50:11.459520
UCHL1|ubiquitin C-terminal hydrolase L1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459523
VLDLR|very low density lipoprotein receptor|'3|Ataxia - paediatric|2021-10-05 10:50:11.459525
VPS13D|vacuolar protein sorting 13 homolog D|'3|Ataxia - paediatric|2021-10-05 10:50:11.459528
VPS41|VPS41, HOPS complex subunit|'3|Ataxia - paediatric|2021-10-05 10:50:11.459530
WDR73|WD repeat domain 73|'3|Ataxia - paediatric|2021-10-05 10:50:11.459534
WDR81|WD repeat domain 81|'3|Ataxia - paediatric|2021-10-05 10:50:11.459537
WFS1|wolframin ER transmembrane glycoprotein|'3|Ataxia - paediatric|2021-10-05 10:50:11.459540
WWOX|WW domain containing oxidoreductase|'3|Ataxia - paediatric|2021-10-05 10:50:11.459542
XRCC1|X-ray repair cross complementing 1|'3|Ataxia - paediatric|2021-10-05 10:50:11.459545
ACBD5|acyl-CoA binding domain containing 5|'2|Ataxia - paediatric|2021-10-05 10:50:11.459547
ATG5|autophagy related 5|'2|Ataxia - paediatric|2021-10-05 10:50:11.459550
ATP2B3|ATPase plasma membrane Ca2+ transporting 3|'2|Ataxia - paediatric|2021-10-05 10:50:11.459552
HARS|histidyl-tRNA synthetase|'2|Ataxia - paediatric|2021-10-05 10:50:11.459555
KCNQ2|potassium voltage-gated channel subfamily Q member 2|'2|Ataxia - paediatric|2021-10-05 10:50:11.459557
MAPK8IP3|mitogen-activated protein kinase 8 interacting protein 3|'2|Ataxia - paediatric|2021-10-05 10:50:11.459559
MKKS|McKusick-Kaufman syndrome|'2|Ataxia - paediatric|2021-10-05 10:50:11.459562
MTPAP|mitochondrial poly(A) polymerase|'2|Ataxia - paediatric|2021-10-05 10:50:11.459564
RFXANK|regulatory factor X associated ankyrin containing protein|'2|Ataxia - paediatric|2021-10-05 10:50:11.459567
SLC9A1|solute carrier family 9 member A1|'2|Ataxia - paediatric|2021-10-05 10:50:11.459569
SYNGAP1|synaptic Ras GTPase activating protein 1|'2|Ataxia - paediatric|2021-10-05 10:50:11.459572
TMEM138|transmembrane protein 138|'2|Ataxia - paediatric|2021-10-05 10:50:11.459574
TMEM231|transmembrane protein 231|'2|Ataxia - paediatric|2021-10-05 10:50:11.459577
UBA5|ubiquitin like modifier activating enzyme 5|'2|Ataxia - paediatric|2021-10-05 10:50:11.459579
UBR4|ubiquitin protein ligase E3 component n-recognin 4|'2|Ataxia - paediatric|2021-10-05 10:50:11.459582
VRK1|vaccinia related kinase 1|'2|Ataxia - paediatric|2021-10-05 10:50:11.459584
AMPD2|adenosine monophosphate deaminase 2|'1|Ataxia - paediatric|2021-10-05 10:50:11.459587
ARL6|ADP ribosylation factor like GTPase 6|'1|Ataxia - paediatric|2021-10-05 10:50:11.459589
ARMC9|armadillo repeat containing 9|'1|Ataxia - paediatric|2021-10-05 10:50:11.459591
BBS10|Bardet-Biedl syndrome 10|'1|Ataxia - paediatric|2021-10-05 10:50:11.459594
BBS12|Bardet-Biedl syndrome 12|'1|Ataxia - paediatric|2021-10-05 10:50:11.459596
BBS2|Bardet-Biedl syndrome 2|'1|Ataxia - paediatric|2021-10-05 10:50:11.459599
BBS4|Bardet-Biedl syndrome 4|'1|Ataxia - paediatric|2021-10-05 10:50:11.459601
BBS5|Bardet-Biedl syndrome 5|'1|Ataxia - paediatric|2021-10-05 10:50:11.459604
BBS7|Bardet-Biedl syndrome 7|'1|Ataxia - paediatric|2021-10-05 10:50:11.459607
BBS9|Bardet-Biedl syndrome 9|'1|Ataxia - paediatric|2021-10-05 10:50:11.459609
CCDC28B|coiled-coil domain containing 28B|'1|Ataxia - paediatric|2021-10-05 10:50:11.459612
CHMP1A|charged multivesicular body protein 1A|'1|Ataxia - paediatric|2021-10-05 10:50:11.459614
COQ5|coenzyme Q5, methyltransferase|'1|Ataxia - paediatric|2021-10-05 10:50:11.459617
CYP2U1|cytochrome P450 family 2 subfamily U member 1|'1|Ataxia - paediatric|2021-10-05 10:50:11.459619
DAB1|DAB1, reelin adaptor protein|'1|Ataxia - paediatric|2021-10-05 10:50:11.459622
ELOVL1|ELOVL fatty acid elongase 1|'1|Ataxia - paediatric|2021-10-05 10:50:11.459624
EXOSC3|exosome component 3|'1|Ataxia - paediatric|2021-10-05 10:50:11.459626
HARS2|histidyl-tRNA synthetase 2, mitochondrial|'1|Ataxia - paediatric|2021-10-05 10:50:11.459629
PCDH12|protocadherin 12|'1|Ataxia - paediatric|2021-10-05 10:50:11.459632
PCYT2|phosphate cytidylyltransferase 2, ethanolamine|'1|Ataxia - paediatric|2021-10-05 10:50:11.459634
RARS2|arginyl-tRNA synthetase 2, mitochondrial|'1|Ataxia - paediatric|2021-10-05 10:50:11.459637
SAR1B|secretion associated Ras related GTPase 1B|'1|Ataxia - paediatric|2021-10-05 10:50:11.459639
TRIM32|tripartite motif containing 32|'1|Ataxia - paediatric|2021-10-05 10:50:11.459642
TSEN2|tRNA splicing endonuclease subunit 2|'1|Ataxia - paediatric|2021-10-05 10:50:11.459644
TSEN34|tRNA splicing endonuclease subunit 34|'1|Ataxia - paediatric|2021-10-05 10:50:11.459647
TTC8|tetratricopeptide repeat domain 8|'1|Ataxia - paediatric|2021-10-05 10:50:11.459649
TTI1|TELO2 interacting protein 1|'1|Ataxia - paediatric|2021-10-05 10:50:11.459652
TUBA1A|tubulin alpha 1a|'1|Ataxia - paediatric|2021-10-05 10:50:11.459654
TUBB2A|tubulin beta 2A class IIa|'1|Ataxia - paediatric|2021-10-05 10:50:11.459657
WDPCP|WD repeat containing planar cell polarity effector|'1|Ataxia - paediatric|2021-10-05 10:50:11.459659
ZNF423|zinc finger protein 423|'1|Ataxia - paediatric|2021-10-05 10:50:11.459661
ZNF592|zinc finger protein 592|'1|Ataxia - paediatric|2021-10-05 10:50:11.459664
'@
$lines = $text -split "`r`n|`n"
$n = $lines.Length
Write-Host "Row count:" $n

$bcdef = New-Object 'object[,]' $n,5
$acol = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $parts = $lines[$i] -split '\|'
    $bcdef[$i,0] = $parts[0]
    $bcdef[$i,1] = $parts[1]
    $bcdef[$i,2] = $parts[2]
    $bcdef[$i,3] = $parts[3]
    $bcdef[$i,4] = $parts[4]
    $acol[$i,0] = $i
}

$startRow = 2
$endRow = $startRow + $n - 1

$rngBF = $ws.Range("B" + $startRow + ":F" + $endRow)
$rngBF.Value = $bcdef
$rngBF.Style = "Normal"

$rngA = $ws.Range("A" + $startRow + ":A" + $endRow)
$rngA.Value = $acol

# ---------------------------------------------------------------
# 3. Apply the header/index style (bold, bordered, centered) to the
#    "A" column cells for the two brand-new rows (225 and 226) by
#    copying the format already used by the rest of column A.
# ---------------------------------------------------------------
$ws.Range("A224").Copy()
$ws.Range("A225:A226").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

Write-Host "Done. Dimension rows:" $ws.UsedRange.Rows.Count "cols:" $ws.UsedRange.Columns.Count
